$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Remove the "Set Assignee" rule (row 21) entirely - entire row delete shifts
# everything below it up by one.
$ws.Rows.Item(21).Delete()

$ws.Range("C22").Select()
